# Update trap card effect texts (Sheet1, column D) to match the new
# "same room" (同房间) wording, and move the active selection from D16 to D15
# (also scrolling the view back to the top, i.e. no topLeftCell override).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value  = "开战时，触发：同房间所有怪物牌点数变为1。"
$ws.Range("D4").Value  = "回合结束时，触发：如果同房间怪物牌合计数量不小于本牌点数，消灭那些怪物牌，然后消灭本牌。"
$ws.Range("D5").Value  = "回合结束时，触发：横置本牌，然后同房间内所有其他牌点数减1。"
$ws.Range("D6").Value  = "回合结束时，触发：横置同房间内所有牌。"
$ws.Range("D7").Value  = "同房间内怪物牌移入或移出时，触发：那张怪物牌点数减2，然后本牌点数减1。"
$ws.Range("D8").Value  = "死亡时，触发：如果本牌点数大于1，则将本牌移动到相邻房间中而不是送墓，然后使本牌和目标房间中所有其他牌点数减1。"
$ws.Range("D9").Value  = "开战时，触发：如果同房间怪物牌合计数量大于1，则将那些怪物牌移入其他房间的同侧，然后消灭本牌。"
$ws.Range("D10").Value = "点数降低时，触发：同房间所有牌点数减1，然后消灭本牌。"
$ws.Range("D11").Value = "同房间内牌移入或移出时，触发：同房间内所有牌点数减1，然后消灭本牌。"
$ws.Range("D12").Value = "死亡时，触发：选1个相邻房间，使其中的所有牌点数减1。"
$ws.Range("D13").Value = "回合结束时：移动1格，然后点数减1。<br>`n点数为0时，触发：消灭同房间1张其他牌。"
$ws.Range("D14").Value = "回合结束时，触发：选同侧至多1张怪物牌，与本牌一起移动到对侧。"

# Scroll the view back to the top and move the selection to D15 (previously
# the view was scrolled to A11 with D16 selected).
$ws.Range("A1").Select() | Out-Null
$ws.Range("D15").Select() | Out-Null
